$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three trailing rows that no longer exist after the re-run
# (old rows 116-118); this shifts nothing else, just shrinks the used range.
$ws.Rows("116:118").Delete()

# New accuracy values for B2:B115 (row order = epoch 0..99 then placeholder rows)
$values = @(0.703125,0.6875,0.46875,0.375,0.3125,0.296875,0.3125,0.34375,0.3125,0.296875,0.25,0.296875,0.40625,0.390625,0.421875,0.328125,0.3125,0.3125,0.34375,0.3125,0.328125,0.265625,0.28125,0.40625,0.40625,0.40625,0.34375,0.3125,0.34375,0.296875,0.25,0.28125,0.265625,0.28125,0.3125,0.265625,0.25,0.25,0.234375,0.234375,0.25,0.234375,0.234375,0.234375,0.234375,0.25,0.25,0.25,0.265625,0.265625,0.28125,0.28125,0.28125,0.28125,0.28125,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.296875,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.28125,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.265625,0.296875,0.234375,0.359375,0.234375,0.34375,0.296875,0.140625,0.34375,0.25,0.375,0.328125,0.390625,0.4545454545454545)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Refresh the repr text in column A for the placeholder "DisplayOutputs" rows (102-115)
# to reflect the new object memory address from the re-run.
for ($row = 102; $row -le 115; $row++) {
    $ws.Cells.Item($row, 1).Value = "<__main__.DisplayOutputs object at 0x7f9304f95b80>"
}
